$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new "Number of Requests" column in H, duplicating the data
# that used to live in column G, then turn G into the rounded version.
# ------------------------------------------------------------------

# 1) H36 gets the original header text ("Number of Requests"), formatted
#    like the rest of the header row (copy format from G36).
$ws.Cells.Item(36, 8).Value = "Number of Requests"
$ws.Cells.Item(36, 7).Copy()
$ws.Cells.Item(36, 8).PasteSpecial(-4122)

# 2) H37:H44 get the same numbers that are currently in G37:G44.
for ($r = 37; $r -le 44; $r++) {
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 7).Value2
}
$ws.Range("G37:G44").Copy()
$ws.Range("H37:H44").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) G36 becomes the new "Number of Requests(Rounded)" header.
$ws.Cells.Item(36, 7).Value = "Number of Requests(Rounded)"

# 4) Column G (values) and the blank column D alongside it now display
#    rounded (integer) numbers.
$ws.Range("D37:D44").NumberFormat = "0"
$ws.Range("G37:G44").NumberFormat = "0"

# ------------------------------------------------------------------
# Cosmetic bits: column widths / view / page setup
# ------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 28.43
$ws.Columns.Item(8).ColumnWidth = 21.17

$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Range("D46").Select() | Out-Null

$ws.PageSetup.Orientation = 1
